$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.548.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.524.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.563"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.96%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.517"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0800"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.912.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.502.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.808"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.549.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("E19").Value = "  -2.60%  "
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0939"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.01%  "
$ws.Range("E24").Value = "  -2.23%  "
$ws.Range("E25").Value = "  -2.77%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.12%  "
$ws.Range("E28").Value = "  -4.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("E37").Value = "  -4.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.60"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("E39").Value = "  -2.71%  "
$ws.Range("E40").Value = "  -0.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.003.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.89%  "
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.763.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "79.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.20%  "
$ws.Range("E50").Value = "  -2.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.03%  "
